$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 54
$ws.Range("C2").Value = "house/house000.jpg"
$ws.Range("D2").Value = "heben"
$ws.Range("E2").Value = "house"
$ws.Range("B3").Value = 87
$ws.Range("C3").Value = "dog/dog009.jpg"
$ws.Range("D3").Value = "grenzen"
$ws.Range("E3").Value = "dog"
$ws.Range("B4").Value = 107
$ws.Range("C4").Value = "dog/dog004.jpg"
$ws.Range("D4").Value = "sparen"
$ws.Range("E4").Value = "dog"
$ws.Range("B5").Value = 9
$ws.Range("C5").Value = "dog/dog027.jpg"
$ws.Range("D5").Value = "fließen"
$ws.Range("E5").Value = "dog"
$ws.Range("B6").Value = 51
$ws.Range("C6").Value = "dog/dog030.jpg"
$ws.Range("D6").Value = "legen"
$ws.Range("E6").Value = "dog"
$ws.Range("B7").Value = 36
$ws.Range("C7").Value = "dog/dog003.jpg"
$ws.Range("D7").Value = "hören"
$ws.Range("E7").Value = "dog"
$ws.Range("B8").Value = 100
$ws.Range("C8").Value = "dog/dog015.jpg"
$ws.Range("D8").Value = "rufen"
$ws.Range("E8").Value = "dog"
$ws.Range("B9").Value = 126
$ws.Range("C9").Value = "house/house016.jpg"
$ws.Range("D9").Value = "hacken"
$ws.Range("E9").Value = "house"
$ws.Range("B10").Value = 76
$ws.Range("C10").Value = "dog/dog017.jpg"
$ws.Range("D10").Value = "piepen"
$ws.Range("E10").Value = "dog"
$ws.Range("B11").Value = 27
$ws.Range("C11").Value = "dog/dog000.jpg"
$ws.Range("D11").Value = "achten"
$ws.Range("E11").Value = "dog"
$ws.Range("B12").Value = 11
$ws.Range("C12").Value = "house/house001.jpg"
$ws.Range("D12").Value = "süßen"
$ws.Range("E12").Value = "house"
$ws.Range("B13").Value = 18
$ws.Range("C13").Value = "house/house027.jpg"
$ws.Range("D13").Value = "zielen"
$ws.Range("E13").Value = "house"
$ws.Range("B14").Value = 30
$ws.Range("C14").Value = "dog/dog002.jpg"
$ws.Range("D14").Value = "kennen"
$ws.Range("E14").Value = "dog"
$ws.Range("B15").Value = 67
$ws.Range("C15").Value = "house/house030.jpg"
$ws.Range("D15").Value = "wählen"
$ws.Range("E15").Value = "house"
$ws.Range("B16").Value = 103
$ws.Range("C16").Value = "house/house024.jpg"
$ws.Range("D16").Value = "rechnen"
$ws.Range("E16").Value = "house"
$ws.Range("B17").Value = 109
$ws.Range("C17").Value = "dog/dog026.jpg"
$ws.Range("D17").Value = "danken"
$ws.Range("E17").Value = "dog"
$ws.Range("B18").Value = 72
$ws.Range("C18").Value = "house/house020.jpg"
$ws.Range("D18").Value = "binden"
$ws.Range("E18").Value = "house"
$ws.Range("B19").Value = 22
$ws.Range("C19").Value = "house/house003.jpg"
$ws.Range("D19").Value = "duschen"
$ws.Range("E19").Value = "house"
$ws.Range("B20").Value = 121
$ws.Range("C20").Value = "dog/dog020.jpg"
$ws.Range("D20").Value = "mögen"
$ws.Range("E20").Value = "dog"
$ws.Range("B21").Value = 94
$ws.Range("C21").Value = "house/house004.jpg"
$ws.Range("D21").Value = "deuten"
$ws.Range("E21").Value = "house"
$ws.Range("B22").Value = 40
$ws.Range("C22").Value = "house/house026.jpg"
$ws.Range("D22").Value = "münzen"
$ws.Range("E22").Value = "house"
$ws.Range("B23").Value = 117
$ws.Range("C23").Value = "house/house014.jpg"
$ws.Range("D23").Value = "streifen"
$ws.Range("E23").Value = "house"
$ws.Range("B24").Value = 3
$ws.Range("C24").Value = "dog/dog028.jpg"
$ws.Range("D24").Value = "leeren"
$ws.Range("E24").Value = "dog"
$ws.Range("B25").Value = 89
$ws.Range("C25").Value = "house/house011.jpg"
$ws.Range("D25").Value = "öffnen"
$ws.Range("E25").Value = "house"
$ws.Range("B26").Value = 115
$ws.Range("C26").Value = "dog/dog012.jpg"
$ws.Range("D26").Value = "lächeln"
$ws.Range("E26").Value = "dog"
$ws.Range("B27").Value = 59
$ws.Range("C27").Value = "dog/dog016.jpg"
$ws.Range("D27").Value = "wachsen"
$ws.Range("E27").Value = "dog"
$ws.Range("B28").Value = 44
$ws.Range("C28").Value = "house/house025.jpg"
$ws.Range("D28").Value = "holen"
$ws.Range("E28").Value = "house"
$ws.Range("B29").Value = 90
$ws.Range("C29").Value = "house/house009.jpg"
$ws.Range("D29").Value = "bergen"
$ws.Range("E29").Value = "house"
$ws.Range("B30").Value = 13
$ws.Range("C30").Value = "house/house031.jpg"
$ws.Range("D30").Value = "trotzen"
$ws.Range("E30").Value = "house"
$ws.Range("B31").Value = 43
$ws.Range("C31").Value = "dog/dog031.jpg"
$ws.Range("D31").Value = "wehen"
$ws.Range("E31").Value = "dog"
$ws.Range("B32").Value = 62
$ws.Range("C32").Value = "house/house013.jpg"
$ws.Range("D32").Value = "stoppen"
$ws.Range("E32").Value = "house"
$ws.Range("B33").Value = 6
$ws.Range("C33").Value = "dog/dog014.jpg"
$ws.Range("D33").Value = "kranken"
$ws.Range("E33").Value = "dog"
